$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The checklist header in column C was "Required"; rename it to "Priority".
$ws.Range("C1").Value = "Priority"

# Move/restore the active selection to the header cell C1 (was sitting on C2).
$ws.Range("C1").Select()
